# Updated cryptos list on Sat Mar 16 02:58:47 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.128.17"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "3.733.98"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'616.60"
$ws.Range("E5").Value = "  +5.98%  "

$ws.Range("D6").Value = "'186.05"
$ws.Range("E6").Value = "  +4.30%  "

$ws.Range("D7").Value = "3.733.48"
$ws.Range("E7").Value = "  -5.18%  "

$ws.Range("D8").Value = "'0.641"
$ws.Range("E8").Value = "  -0.75%  "

$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("D10").Value = "'0.721"
$ws.Range("E10").Value = "  -1.17%  "

$ws.Range("E11").Value = "  -4.81%  "

$ws.Range("D12").Value = "'56.84"
$ws.Range("E12").Value = "  +6.02%  "

$ws.Range("E13").Value = "  -4.76%  "

$ws.Range("D14").Value = "'10.68"
$ws.Range("E14").Value = "  -2.09%  "

$ws.Range("D15").Value = "4.327.64"
$ws.Range("E15").Value = "  -0.32%  "

$ws.Range("D16").Value = "3.733.53"
$ws.Range("E16").Value = "  -2.41%  "

$ws.Range("E17").Value = "  -1.64%  "

$ws.Range("D18").Value = "'13.09"
$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("E19").Value = "  -0.91%  "

$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("D21").Value = "68.967.82"
$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("D22").Value = "'413.95"
$ws.Range("E22").Value = "  -1.09%  "

$ws.Range("D23").Value = "'4.68"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("D24").Value = "'89.67"
$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("E25").Value = "  -1.87%  "

$ws.Range("D26").Value = "'12.84"
$ws.Range("E26").Value = "  -2.33%  "

$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("D28").Value = "'6.07"
$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("E30").Value = "  -2.30%  "

$ws.Range("E31").Value = "  -1.21%  "

$ws.Range("D32").Value = "'7.31"
$ws.Range("E32").Value = "  -13.23%  "

$ws.Range("D33").Value = "'12.72"
$ws.Range("E33").Value = "  -2.52%  "

$ws.Range("E34").Value = "  +1.86%  "

$ws.Range("D35").Value = "'621.49"
$ws.Range("E35").Value = "  +2.24%  "

$ws.Range("D36").Value = "'44.56"
$ws.Range("E36").Value = "  -2.54%  "

$ws.Range("D37").Value = "'66.00"
$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("D38").Value = "0.0₃0872"
$ws.Range("E38").Value = "  -8.22%  "

# Rows 39/40: Dai and TheGraph swap places (with new price/volume data)
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.406"
$ws.Range("E39").Value = "  -0.94%  "

$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("E42").Value = "  +2.04%  "

$ws.Range("E43").Value = "  -2.22%  "

$ws.Range("D44").Value = "'0.0444"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("D47").Value = "'9.25"
$ws.Range("E47").Value = "  -4.44%  "

$ws.Range("D48").Value = "2.843.74"
$ws.Range("E48").Value = "  +2.36%  "

$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("E50").Value = "  -16.88%  "

$ws.Range("D51").Value = "'3.13"
$ws.Range("E51").Value = "  -2.66%  "
